# Employment_Youth_unemployment_rate_OECD_2013 update:
#  - rename worksheet (truncated sheet name)
#  - fix the metadata key/value misalignment ("indicator source treated as
#    array") so several target$/scoring$/source$ rows carry their correct
#    values
#  - bump the scoring timestamp + data-host year to the new data refresh
#  - give the scoring$timestamp date cell an explicit dd/mm/yyyy format

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (31-char Excel tab-name limit truncates the csv name) ---
$ws.Name = "Employment_Youth_unemployment_r"

# --- Metadata column B fixes (column A labels are untouched) ---

# target$explanation$de / target$explanation$en
$ws.Range("B17").Value = "Der Zielwert orientiert sich an den Ländern, die aktuell am besten abschneiden. "
$ws.Range("B18").Value = "The target value is based on the currently best performing countries. "

# target$ministerial_responsibility
$ws.Range("B22").Value = "BMAS"

# scoring$timestamp (now an actual date, formatted dd/mm/yyyy)
$ws.Range("B25").Value = 42711.0
$ws.Range("B25").NumberFormat = "dd/mm/yyyy"

# scoring$type
$ws.Range("B26").Value = "national"

# scoring$timestamp_data_host
$ws.Range("B27").Value = 2015.0

# source (no value directly on the parent key)
$ws.Range("B28").Value = ""

# source$type
$ws.Range("B29").Value = "inofficial"

# source$maintainer
$ws.Range("B34").Value = "OKF"

# source$license
$ws.Range("B35").Value = "The OECD supports free use and consultation of its data by the public. Information source must be cited. Download is feasible. "
